$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy header style (bold/border/centered) from an existing header cell (F1) to G1:H1
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Updated metric values in row 2
$ws.Range("B2").Value = 0.3578863750850889
$ws.Range("C2").Value = 0.997868172976735
$ws.Range("D2").Value = 0.4984859599743723

# Updated model description text
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=100))])"

# New data cells for elapsed time / cpu
$ws.Range("G2").Value = 0.1256850772835605
$ws.Range("H2").Value = 0.99
